# AB - Added first 4 zones of row 5
#
# The canonical diff for this commit only touches the cached text of the
# "datetimeFigureOut" date fields that live on the Slide Master, every
# Slide Layout, and the Notes Master (PowerPoint re-stamps these cached
# field values whenever the deck is saved on a different day - the date
# simply rolled from 21/10/2016 to 22/10/2016). Walk every one of those
# containers, find the "Date" placeholder shape on each, and refresh its
# text to the new cached date.

$p = $ppt.ActivePresentation

$oldDate = "21/10/2016"
$newDate = "22/10/2016"

function Update-DatePlaceholder($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $shp = $container.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 16) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# 1) Slide Master
Update-DatePlaceholder($p.SlideMaster)

# 2) Every Slide Layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder($layouts.Item($i))
}

# 3) Notes Master
Update-DatePlaceholder($p.NotesMaster)
